# Auto-generated edit script: updates cached profit-calculation values
# in the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (per-recipe cost/profit
# columns H:N), matching the scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1460.8064
$ws.Range("I80").Value = 864.9091
$ws.Range("J80").Value = 1788.55
$ws.Range("K80").Value = 2594.7273
$ws.Range("L80").Value = 5365.65
$ws.Range("M80").Value = -1596.7273
$ws.Range("N80").Value = -7361.65
$ws.Range("H83").Value = 1460.8064
$ws.Range("I83").Value = 864.9091
$ws.Range("J83").Value = 1788.55
$ws.Range("K83").Value = 7784.1819
$ws.Range("L83").Value = 16096.95
$ws.Range("M83").Value = -2792.1819
$ws.Range("N83").Value = -26080.95
$ws.Range("H100").Value = 3319.8667
$ws.Range("I100").Value = 2454.9092
$ws.Range("J100").Value = 5698.5
$ws.Range("K100").Value = 2454.9092
$ws.Range("L100").Value = 5698.5
$ws.Range("M100").Value = -1913.9092
$ws.Range("N100").Value = -6780.5
$ws.Range("H137").Value = 1106.0869
$ws.Range("I137").Value = 1051.2941
$ws.Range("K137").Value = 3153.8823
$ws.Range("M137").Value = -603.8823000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3158.6924
$ws.Range("I2").Value = 2778.9092
$ws.Range("J2").Value = 5247.5
$ws.Range("K2").Value = 2778.9092
$ws.Range("L2").Value = 5247.5
$ws.Range("M2").Value = -2665.9092
$ws.Range("N2").Value = -5473.5
$ws.Range("H32").Value = 10535.553
$ws.Range("I32").Value = 4009.75
$ws.Range("K32").Value = 4009.75
$ws.Range("M32").Value = -3722.75
$ws.Range("H45").Value = 9467.666999999999
$ws.Range("I45").Value = 12940.556
$ws.Range("K45").Value = 12940.556
$ws.Range("M45").Value = -12563.556
$ws.Range("H54").Value = 39495
$ws.Range("J54").Value = 39495
$ws.Range("L54").Value = 39495
$ws.Range("N54").Value = -41033
$ws.Range("H97").Value = 2445.8
$ws.Range("I97").Value = 2120.4092
$ws.Range("K97").Value = 2120.4092
$ws.Range("M97").Value = -1624.4092
$ws.Range("H102").Value = 1841.4828
$ws.Range("I102").Value = 1818.7307
$ws.Range("K102").Value = 1818.7307
$ws.Range("M102").Value = -196.7307000000001
$ws.Range("H110").Value = 950.35
$ws.Range("I110").Value = 989.2778
$ws.Range("J110").Value = 600
$ws.Range("K110").Value = 989.2778
$ws.Range("L110").Value = 600
$ws.Range("M110").Value = 1055.7222
$ws.Range("N110").Value = -4690
$ws.Range("H116").Value = 3158.6924
$ws.Range("I116").Value = 2778.9092
$ws.Range("J116").Value = 5247.5
$ws.Range("K116").Value = 2778.9092
$ws.Range("L116").Value = 5247.5
$ws.Range("M116").Value = -484.9092000000001
$ws.Range("N116").Value = -9835.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3158.6924
$ws.Range("I3").Value = 2778.9092
$ws.Range("J3").Value = 5247.5
$ws.Range("K3").Value = 2778.9092
$ws.Range("L3").Value = 5247.5
$ws.Range("M3").Value = -2664.9092
$ws.Range("N3").Value = -5475.5
$ws.Range("H23").Value = 2651
$ws.Range("I23").Value = 312
$ws.Range("J23").Value = 4990
$ws.Range("K23").Value = 312
$ws.Range("L23").Value = 4990
$ws.Range("M23").Value = -29
$ws.Range("N23").Value = -5556
$ws.Range("H86").Value = 2907.5
$ws.Range("I86").Value = 2862
$ws.Range("K86").Value = 2862
$ws.Range("M86").Value = -1739
$ws.Range("H89").Value = 2907.5
$ws.Range("I89").Value = 2862
$ws.Range("K89").Value = 14310
$ws.Range("M89").Value = -8694
$ws.Range("H99").Value = 900
$ws.Range("I99").Value = 900
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 900
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 598
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 50010
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 1565.2693
$ws.Range("I134").Value = 1340.15
$ws.Range("J134").Value = 2315.6667
$ws.Range("K134").Value = 4020.45
$ws.Range("L134").Value = 6947.000100000001
$ws.Range("M134").Value = -1485.45
$ws.Range("N134").Value = -12017.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4559.6
$ws.Range("I16").Value = 5999.6665
$ws.Range("K16").Value = 5999.6665
$ws.Range("M16").Value = -5712.6665
$ws.Range("H22").Value = 476.75
$ws.Range("I22").Value = 283
$ws.Range("J22").Value = 748
$ws.Range("K22").Value = 283
$ws.Range("L22").Value = 748
$ws.Range("M22").Value = 67
$ws.Range("N22").Value = -1448
$ws.Range("H62").Value = 3300.889
$ws.Range("I62").Value = 3213.5
$ws.Range("K62").Value = 3213.5
$ws.Range("M62").Value = -2589.5
$ws.Range("H65").Value = 3300.889
$ws.Range("I65").Value = 3213.5
$ws.Range("K65").Value = 16067.5
$ws.Range("M65").Value = -12947.5
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("H107").Value = 1631.52
$ws.Range("I107").Value = 2088.0908
$ws.Range("J107").Value = 1272.7858
$ws.Range("K107").Value = 2088.0908
$ws.Range("L107").Value = 1272.7858
$ws.Range("M107").Value = -168.0907999999999
$ws.Range("N107").Value = -5112.7858
$ws.Range("H113").Value = 4559.6
$ws.Range("I113").Value = 5999.6665
$ws.Range("K113").Value = 5999.6665
$ws.Range("M113").Value = -3829.6665
$ws.Range("H122").Value = 70313
$ws.Range("I122").Value = 112066.555
$ws.Range("K122").Value = 336199.665
$ws.Range("M122").Value = -333749.665
$ws.Range("H132").Value = 3895.5
$ws.Range("I132").Value = 3943.1304
$ws.Range("J132").Value = 3676.4
$ws.Range("K132").Value = 11829.3912
$ws.Range("L132").Value = 11029.2
$ws.Range("M132").Value = -9299.3912
$ws.Range("N132").Value = -16089.2
$ws.Range("H134").Value = 2593.3215
$ws.Range("I134").Value = 1622.1
$ws.Range("K134").Value = 4866.299999999999
$ws.Range("M134").Value = -2331.299999999999
$ws.Range("H141").Value = 268864.7
$ws.Range("J141").Value = 268864.7
$ws.Range("L141").Value = 268864.7
$ws.Range("N141").Value = -279224.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1686.75
$ws.Range("J5").Value = 4212.6665
$ws.Range("L5").Value = 12637.9995
$ws.Range("N5").Value = -12861.9995
$ws.Range("H86").Value = 655.26086
$ws.Range("I86").Value = 467.36365
$ws.Range("J86").Value = 827.5
$ws.Range("K86").Value = 1402.09095
$ws.Range("L86").Value = 2482.5
$ws.Range("M86").Value = -216.09095
$ws.Range("N86").Value = -4854.5
$ws.Range("H89").Value = 655.26086
$ws.Range("I89").Value = 467.36365
$ws.Range("J89").Value = 827.5
$ws.Range("K89").Value = 4206.27285
$ws.Range("L89").Value = 7447.5
$ws.Range("M89").Value = 1721.72715
$ws.Range("N89").Value = -19303.5
$ws.Range("H135").Value = 1686.75
$ws.Range("J135").Value = 4212.6665
$ws.Range("L135").Value = 37913.9985
$ws.Range("N135").Value = -42983.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -38744
$ws.Range("H122").Value = 4012.5
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 4275
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 12825
$ws.Range("M122").Value = -8800
$ws.Range("N122").Value = -17725
$ws.Range("H132").Value = 7170.7144
$ws.Range("I132").Value = 7170.7144
$ws.Range("K132").Value = 21512.1432
$ws.Range("M132").Value = -18982.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 29622.357
$ws.Range("I100").Value = 5692.125
$ws.Range("K100").Value = 5692.125
$ws.Range("M100").Value = -5151.125
$ws.Range("H132").Value = 2977.8438
$ws.Range("I132").Value = 2783.8518
$ws.Range("K132").Value = 8351.555399999999
$ws.Range("M132").Value = -5821.555399999999
$ws.Range("H136").Value = 4347.4
$ws.Range("I136").Value = 3841.2666
$ws.Range("J136").Value = 5865.8
$ws.Range("K136").Value = 11523.7998
$ws.Range("L136").Value = 17597.4
$ws.Range("M136").Value = -8973.799800000001
$ws.Range("N136").Value = -22697.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3540.1
$ws.Range("I100").Value = 3782.182
$ws.Range("J100").Value = 3244.2222
$ws.Range("K100").Value = 7564.364
$ws.Range("L100").Value = 6488.4444
$ws.Range("M100").Value = -7023.364
$ws.Range("N100").Value = -7570.4444
$ws.Range("H107").Value = 22728436
$ws.Range("J107").Value = 50001132
$ws.Range("L107").Value = 150003396
$ws.Range("N107").Value = -150007236
$ws.Range("H136").Value = 3349.5862
$ws.Range("I136").Value = 3101.9524
$ws.Range("J136").Value = 3999.625
$ws.Range("K136").Value = 9305.8572
$ws.Range("L136").Value = 11998.875
$ws.Range("M136").Value = -6755.8572
$ws.Range("N136").Value = -17098.875
